$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 104, shifting the existing
# rows 104-110 down to 105-111 (all their data stays intact).
$ws.Rows(104).Insert()

# Populate the newly inserted row 104 with the new weekly record.
$ws.Cells.Item(104, 1).Value = 10
$ws.Cells.Item(104, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(104, 3).Value = "La Araucanía"
$ws.Cells.Item(104, 4).Value = 44783
$ws.Cells.Item(104, 5).Value = 9
$ws.Cells.Item(104, 6).Value = 100112035
$ws.Cells.Item(104, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(104, 8).Value = "Sin especificar"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 50
$ws.Cells.Item(104, 11).Value = 25000
$ws.Cells.Item(104, 12).Value = 26000
$ws.Cells.Item(104, 13).Value = 25400
$ws.Cells.Item(104, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(104, 15).Value = "Región Metropolitana"
$ws.Cells.Item(104, 16).Value = 2540
$ws.Cells.Item(104, 17).Value = 10
$ws.Cells.Item(104, 18).Value = "Hortaliza"
